# "orden de las listas" - reorder party columns and duplicate them (D:H) into (I:M),
# shifting the trailing "Votos ..." columns from I:L to N:Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers (A1:C1 are untouched by this change) ---
$ws.Range("D1").Value = "CONSENSO FEDERAL"
$ws.Range("E1").Value = "FRENTE DE IZQUIERDA Y DE TRABAJADORES - UNIDAD"
$ws.Range("F1").Value = "FRENTE DE TODOS"
$ws.Range("G1").Value = "JUNTOS POR EL CAMBIO"
$ws.Range("H1").Value = "UNITE POR LA LIBERTAD Y LA DIGNIDAD"

$ws.Range("I1").Value = "CONSENSO FEDERAL"
$ws.Range("J1").Value = "FRENTE DE IZQUIERDA Y DE TRABAJADORES - UNIDAD"
$ws.Range("K1").Value = "FRENTE DE TODOS"
$ws.Range("L1").Value = "JUNTOS POR EL CAMBIO"
$ws.Range("M1").Value = "UNITE POR LA LIBERTAD Y LA DIGNIDAD"

$ws.Range("N1").Value = "Votos Nulos"
$ws.Range("O1").Value = "Votos Recurridos"
$ws.Range("P1").Value = "Votos impugnados"
$ws.Range("Q1").Value = "Votos en blanco"

# --- Row 2: values (A2:C2 are untouched by this change) ---
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 110
$ws.Range("G2").Value = 111
$ws.Range("H2").Value = 4

$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 110
$ws.Range("L2").Value = 111
$ws.Range("M2").Value = 4

$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 13
